$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 10.14708026778702
$ws.Range("C2").Value = 5.375521853503479
$ws.Range("D2").Value = 5.992139340925059
$ws.Range("E2").Value = 16.39620419057874
$ws.Range("G2").Value = 3.64752412408999
$ws.Range("K2").Value = 9.392611218413609
$ws.Range("O2").Value = 24.95410262116516
$ws.Range("B3").Value = 9.85068229757689
$ws.Range("C3").Value = 5.174292670948209
$ws.Range("D3").Value = 5.874898842391129
$ws.Range("E3").Value = 15.47245837382427
$ws.Range("G3").Value = 3.64996292289891
$ws.Range("K3").Value = 9.188900869594008
$ws.Range("O3").Value = 24.9647474462114
$ws.Range("B4").Value = 9.666592409111445
$ws.Range("C4").Value = 5.045510473695003
$ws.Range("D4").Value = 5.803535966990759
$ws.Range("E4").Value = 14.88122209026717
$ws.Range("G4").Value = 3.651537378023464
$ws.Range("K4").Value = 9.06370708904592
$ws.Range("O4").Value = 24.97740567019552
$ws.Range("B5").Value = 9.591172527657672
$ws.Range("C5").Value = 4.991759221802056
$ws.Range("D5").Value = 5.774656154355168
$ws.Range("E5").Value = 14.63451316271008
$ws.Range("G5").Value = 3.652198418473267
$ws.Range("K5").Value = 9.012738037479922
$ws.Range("O5").Value = 24.98409925333344
$ws.Range("B6").Value = 9.578628687534383
$ws.Range("C6").Value = 4.982758433311137
$ws.Range("D6").Value = 5.769874135316933
$ws.Range("E6").Value = 14.59320753499653
$ws.Range("G6").Value = 3.65230935972895
$ws.Range("K6").Value = 9.004279887138148
$ws.Range("O6").Value = 24.98530330361847
$ws.Range("B7").Value = 9.665576724450492
$ws.Range("C7").Value = 5.044790652991947
$ws.Range("D7").Value = 5.803145610414774
$ws.Range("E7").Value = 14.87791787206011
$ws.Range("G7").Value = 3.651546214254532
$ws.Range("K7").Value = 9.063019402951324
$ws.Range("O7").Value = 24.97748973231677
$ws.Range("B8").Value = 10.04539164156706
$ws.Range("C8").Value = 5.307246857296024
$ws.Range("D8").Value = 5.951611979748766
$ws.Range("E8").Value = 16.08283683427257
$ws.Range("G8").Value = 3.648349077565916
$ws.Range("K8").Value = 9.322441442477736
$ws.Range("O8").Value = 24.95649996800778
$ws.Range("B9").Value = 10.76838834984202
$ws.Range("C9").Value = 5.77883630259964
$ws.Range("D9").Value = 6.245776082497531
$ws.Range("E9").Value = 18.29804626260722
$ws.Range("G9").Value = 3.642687480382768
$ws.Range("K9").Value = 9.826977902096964
$ws.Range("O9").Value = 24.96407271876498
$ws.Range("B10").Value = 11.27996572198052
$ws.Range("C10").Value = 6.097221164779831
$ws.Range("D10").Value = 6.461235412884745
$ws.Range("E10").Value = 19.92712090215458
$ws.Range("G10").Value = 3.638894109335553
$ws.Range("K10").Value = 10.19092405660235
$ws.Range("O10").Value = 24.99951028970537
$ws.Range("B11").Value = 11.50727902053532
$ws.Range("C11").Value = 6.235661793129106
$ws.Range("D11").Value = 6.558607789958288
$ws.Range("E11").Value = 20.62650663717255
$ws.Range("G11").Value = 3.637246977445189
$ws.Range("K11").Value = 10.35421465036082
$ws.Range("O11").Value = 25.02213562230112
$ws.Range("B12").Value = 11.59249455317131
$ws.Range("C12").Value = 6.287145740972165
$ws.Range("D12").Value = 6.595348698633413
$ws.Range("E12").Value = 20.88538772446913
$ws.Range("G12").Value = 3.636634466139545
$ws.Range("K12").Value = 10.41565954694906
$ws.Range("O12").Value = 25.03163835118458
$ws.Range("B13").Value = 11.57418152468549
$ws.Range("C13").Value = 6.276099882484927
$ws.Range("D13").Value = 6.587442362068481
$ws.Range("E13").Value = 20.82989763900306
$ws.Range("G13").Value = 3.636765883350844
$ws.Range("K13").Value = 10.40244458315182
$ws.Range("O13").Value = 25.02955019940674
$ws.Range("B14").Value = 11.51430751216371
$ws.Range("C14").Value = 6.239916366976446
$ws.Range("D14").Value = 6.561633349836975
$ws.Range("E14").Value = 20.6479242176566
$ws.Range("G14").Value = 3.637196361242404
$ws.Range("K14").Value = 10.35927791643057
$ws.Range("O14").Value = 25.02289868396274
$ws.Range("B15").Value = 11.47751815205305
$ws.Range("C15").Value = 6.217629866567568
$ws.Range("D15").Value = 6.545806265877243
$ws.Range("E15").Value = 20.53568492473746
$ws.Range("G15").Value = 3.637461500859825
$ws.Range("K15").Value = 10.33278453082432
$ws.Range("O15").Value = 25.01894617766088
$ws.Range("B16").Value = 11.26499396211646
$ws.Range("C16").Value = 6.088043204258531
$ws.Range("D16").Value = 6.45485555548068
$ws.Range("E16").Value = 19.8805784013876
$ws.Range("G16").Value = 3.639003327148324
$ws.Range("K16").Value = 10.18020144531347
$ws.Range("O16").Value = 24.99816255213302
$ws.Range("B17").Value = 11.1331689923912
$ws.Range("C17").Value = 6.00689226567659
$ws.Range("D17").Value = 6.398867114100971
$ws.Range("E17").Value = 19.46803835211009
$ws.Range("G17").Value = 3.639969244726455
$ws.Range("K17").Value = 10.08596884828773
$ws.Range("O17").Value = 24.98707850890672
$ws.Range("B18").Value = 11.05684251062538
$ws.Range("C18").Value = 5.959615574465095
$ws.Range("D18").Value = 6.366605941270051
$ws.Range("E18").Value = 19.22683031224161
$ws.Range("G18").Value = 3.64053220684624
$ws.Range("K18").Value = 10.03155857368751
$ws.Range("O18").Value = 24.98131558490399
$ws.Range("B19").Value = 11.03091582234516
$ws.Range("C19").Value = 5.943506002067909
$ws.Range("D19").Value = 6.355674127485137
$ws.Range("E19").Value = 19.14448610597393
$ws.Range("G19").Value = 3.640724087702826
$ws.Range("K19").Value = 10.01310208736873
$ws.Range("O19").Value = 24.97946952040224
$ws.Range("B20").Value = 11.14725482091696
$ws.Range("C20").Value = 6.015593266535117
$ws.Range("D20").Value = 6.404833486952059
$ws.Range("E20").Value = 19.51236001781804
$ws.Range("G20").Value = 3.639865656591749
$ws.Range("K20").Value = 10.09602229312119
$ws.Range("O20").Value = 24.9881950503637
$ws.Range("B21").Value = 11.53191802614811
$ws.Range("C21").Value = 6.250570019721767
$ws.Range("D21").Value = 6.569217967787823
$ws.Range("E21").Value = 20.70153573843455
$ws.Range("G21").Value = 3.637069615411999
$ws.Range("K21").Value = 10.37196807229535
$ws.Range("O21").Value = 25.02482702693205
$ws.Range("B22").Value = 11.77825308749305
$ws.Range("C22").Value = 6.398651290381391
$ws.Range("D22").Value = 6.675867052014972
$ws.Range("E22").Value = 21.44400035218508
$ws.Range("G22").Value = 3.635307620040344
$ws.Range("K22").Value = 10.55001675335283
$ws.Range("O22").Value = 25.05421724488147
$ws.Range("B23").Value = 11.64726899758995
$ws.Range("C23").Value = 6.320125973424819
$ws.Range("D23").Value = 6.619030875095266
$ws.Range("E23").Value = 21.05089979473863
$ws.Range("G23").Value = 3.636242069394264
$ws.Range("K23").Value = 10.4552188656519
$ws.Range("O23").Value = 25.03803289864601
$ws.Range("B24").Value = 11.14088828867717
$ws.Range("C24").Value = 6.011661479697771
$ws.Range("D24").Value = 6.402136313130895
$ws.Range("E24").Value = 19.49233475884239
$ws.Range("G24").Value = 3.639912464972636
$ws.Range("K24").Value = 10.09147785685943
$ws.Range("O24").Value = 24.98768836342697
$ws.Range("B25").Value = 10.57581947844026
$ws.Range("C25").Value = 5.656077057464421
$ws.Range("D25").Value = 6.166137226063226
$ws.Range("E25").Value = 17.68289412145244
$ws.Range("G25").Value = 3.644154463928588
$ws.Range("K25").Value = 9.691380926562539
$ws.Range("O25").Value = 24.95678632520319
